$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.60"
$ws.Range("E2").Value = "'1.28%"
$ws.Range("G2").Value = "'2"

$ws.Range("D3").Value = "'36.17"
$ws.Range("E3").Value = "'-6.63%"
$ws.Range("G3").Value = "'2"

$ws.Range("D4").Value = "'5.040"
$ws.Range("E4").Value = "'1.27%"
$ws.Range("G4").Value = "'2"

$ws.Range("D5").Value = "'0.07826"
$ws.Range("E5").Value = "'1.45%"
$ws.Range("G5").Value = "'2"

$ws.Range("D6").Value = "'2.117"
$ws.Range("E6").Value = "'-3.90%"
$ws.Range("G6").Value = "'2"

$ws.Range("D7").Value = "'7.922"
$ws.Range("E7").Value = "'-1.06%"
$ws.Range("G7").Value = "'2"

$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = "'4.055"
$ws.Range("E8").Value = "'1.29%"
$ws.Range("G8").Value = "'2"

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = "'0.9237"
$ws.Range("E9").Value = "'0.54%"
$ws.Range("G9").Value = "'2"

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.09674"
$ws.Range("E10").Value = "'6.78%"
$ws.Range("G10").Value = "'2"

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1882"
$ws.Range("E11").Value = "'4.70%"
$ws.Range("G11").Value = "'2"

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.08747"
$ws.Range("E12").Value = "'3.71%"
$ws.Range("G12").Value = "'2"

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03495"
$ws.Range("E13").Value = "'-2.52%"
$ws.Range("G13").Value = "'2"

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09922"
$ws.Range("E14").Value = "'-0.27%"
$ws.Range("G14").Value = "'2"

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001442"
$ws.Range("E15").Value = "'-2.50%"
$ws.Range("G15").Value = "'2"

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.005709"
$ws.Range("E16").Value = "'-0.60%"
$ws.Range("G16").Value = "'2"

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.459"
$ws.Range("E17").Value = "'-0.71%"
$ws.Range("G17").Value = "'2"

$ws.Range("D18").Value = "'2.481"
$ws.Range("E18").Value = "'12.27%"
$ws.Range("G18").Value = "'2"

$ws.Range("D19").Value = "'0.3417"
$ws.Range("E19").Value = "'-1.34%"
$ws.Range("G19").Value = "'2"

$ws.Range("D20").Value = "'0.1346"
$ws.Range("E20").Value = "'1.96%"
$ws.Range("G20").Value = "'2"

$ws.Range("E21").Value = "'4.27%"
$ws.Range("G21").Value = "'2"

$ws.Range("D22").Value = "'0.2290"
$ws.Range("E22").Value = "'2.33%"
$ws.Range("G22").Value = "'2"

$ws.Range("D23").Value = "'0.04612"
$ws.Range("E23").Value = "'-1.18%"
$ws.Range("G23").Value = "'2"

$ws.Range("E24").Value = "'14.80%"
$ws.Range("G24").Value = "'2"

$ws.Range("D25").Value = "'0.001229"
$ws.Range("E25").Value = "'-0.28%"
$ws.Range("G25").Value = "'2"

$ws.Range("E26").Value = "'7.44%"
$ws.Range("G26").Value = "'2"

$ws.Range("E27").Value = "'-42.75%"
$ws.Range("G27").Value = "'2"

$ws.Range("G28").Value = "'2"

$ws.Range("G29").Value = "'2"

$ws.Range("G30").Value = "'2"

$ws.Range("G31").Value = "'2"

$ws.Range("G32").Value = "'2"

$ws.Range("G33").Value = "'2"

$ws.Range("G34").Value = "'2"

$ws.Range("G35").Value = "'2"

$ws.Range("G36").Value = "'2"

$ws.Range("G37").Value = "'2"

$ws.Range("G38").Value = "'2"

$ws.Range("D39").Value = "'0.01818"
$ws.Range("E39").Value = "'4.15%"
$ws.Range("G39").Value = "'2"

$ws.Range("D40").Value = "'0.04754"
$ws.Range("E40").Value = "'1.41%"
$ws.Range("G40").Value = "'2"

$ws.Range("D41").Value = "'0.007503"
$ws.Range("E41").Value = "'-5.04%"
$ws.Range("G41").Value = "'2"

$ws.Range("D42").Value = "'0.1401"
$ws.Range("E42").Value = "'0.70%"
$ws.Range("G42").Value = "'2"

$ws.Range("D43").Value = "'0.007719"
$ws.Range("E43").Value = "'3.64%"
$ws.Range("G43").Value = "'2"

$ws.Range("E44").Value = "'-1.98%"
$ws.Range("G44").Value = "'2"

$ws.Range("D45").Value = "'0.01080"
$ws.Range("E45").Value = "'14.04%"
$ws.Range("G45").Value = "'2"

$ws.Range("D46").Value = "'0.00006154"
$ws.Range("E46").Value = "'1.70%"
$ws.Range("G46").Value = "'2"

$ws.Range("E47").Value = "'-0.25%"
$ws.Range("G47").Value = "'2"

$ws.Range("D48").Value = "'0.0005802"
$ws.Range("E48").Value = "'0.03%"
$ws.Range("G48").Value = "'2"

$ws.Range("D49").Value = "'38.56"
$ws.Range("E49").Value = "'635.04%"
$ws.Range("G49").Value = "'2"

$ws.Range("E50").Value = "'-25.85%"
$ws.Range("G50").Value = "'2"

$ws.Range("E51").Value = "'-0.25%"
$ws.Range("G51").Value = "'2"
